# Updated main GSC export data: drop the oldest date's row from the
# "Chart" sheet. Excel will automatically shift the remaining rows up,
# shrink the used range, and drop the now-unreferenced shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
